$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New claim row data (row 3) — a new "Ambiente" (preproducciongestion) entry with
# its own NroPoliza/FechaSiniestro, generated for litigio (mediación o juicio).
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("E3").Value = "'04104014484"
$ws.Range("G3").Value = "'15/03/2021"

# Turn the new URL in B3 into a live hyperlink (matches the other Ambiente/URL rows).
$ws.Range("B3").Hyperlinks.Add($ws.Range("B3"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B3").Style = "Hipervínculo"

# Move the active selection to the new row, as left by the edit session.
$ws.Range("B4").Select()
